# Updates the "cryptos" price table (Sheet1) with a fresh snapshot of
# prices / 1h-volume percentages pulled by the scheduled GitHub Actions
# scraper, and reorders a couple of coins that swapped rank (MXToken /
# HuobiToken and RocketPoolETH / WEMIXToken).
#
# Column D ("Price") and column E ("Volume(1h)") in the source sheet are
# plain text (not real numbers - many prices even contain two thousands
# separators, e.g. "34.175.25"), so every value is written back as text.
# Excel auto-converts a clean decimal-looking string (e.g. "223.43") to a
# number when it is assigned directly, which would both change the cell's
# stored type and silently normalise formatting (e.g. "1.00" -> 1). To
# keep those cells textual - matching the original workbook - a leading
# apostrophe is used to force text entry for any value that looks like a
# plain number, and the resulting "quote prefix" cell format is cleared
# right after so the cell's style stays the same as before the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Test-LooksNumeric($Value) {
    # Matches plain integers/decimals such as "223.43", "0.0722", "1.00".
    # Values with more than one "." (e.g. "34.175.25"), percentages, URLs
    # and coin names never match, so they are written verbatim.
    return ($Value -match '^\s*[+-]?[0-9]+(\.[0-9]+)?\s*$')
}

function Set-CellText($Range, $Value) {
    if (Test-LooksNumeric $Value) {
        $Range.Value = "'" + $Value
        $Range.ClearFormats()
    } else {
        $Range.Value = $Value
    }
}

$updates = @(
    @{ Cell = "D2"; Value = "34.175.25" },
    @{ Cell = "D3"; Value = "1.804.76" },
    @{ Cell = "E3"; Value = "  +0.68%  " },
    @{ Cell = "E4"; Value = "  -0.01%  " },
    @{ Cell = "D5"; Value = "223.43" },
    @{ Cell = "E5"; Value = "  +0.30%  " },
    @{ Cell = "E6"; Value = "  -0.13%  " },
    @{ Cell = "D7"; Value = "1.00" },
    @{ Cell = "E7"; Value = "  -0.02%  " },
    @{ Cell = "D8"; Value = "32.65" },
    @{ Cell = "E8"; Value = "  +1.48%  " },
    @{ Cell = "D9"; Value = "0.288" },
    @{ Cell = "E9"; Value = "  +2.51%  " },
    @{ Cell = "D10"; Value = "0.0722" },
    @{ Cell = "E10"; Value = "  +4.41%  " },
    @{ Cell = "D11"; Value = "0.0929" },
    @{ Cell = "E11"; Value = "  -0.70%  " },
    @{ Cell = "D12"; Value = "2.064.06" },
    @{ Cell = "E12"; Value = "  +0.66%  " },
    @{ Cell = "D13"; Value = "1.806.20" },
    @{ Cell = "E13"; Value = "  +0.88%  " },
    @{ Cell = "D14"; Value = "10.99" },
    @{ Cell = "E14"; Value = "  +1.13%  " },
    @{ Cell = "E15"; Value = "  +0.11%  " },
    @{ Cell = "D16"; Value = "34.204.82" },
    @{ Cell = "E16"; Value = "  -1.25%  " },
    @{ Cell = "E17"; Value = "  -0.98%  " },
    @{ Cell = "D18"; Value = "68.69" },
    @{ Cell = "E18"; Value = "  +0.10%  " },
    @{ Cell = "D19"; Value = "248.04" },
    @{ Cell = "E19"; Value = "  -2.51%  " },
    @{ Cell = "E20"; Value = "  +0.52%  " },
    @{ Cell = "D21"; Value = "11.01" },
    @{ Cell = "E21"; Value = "  +5.38%  " },
    @{ Cell = "E22"; Value = "  -0.05%  " },
    @{ Cell = "E23"; Value = "  -1.00%  " },
    @{ Cell = "D24"; Value = "2.13" },
    @{ Cell = "E24"; Value = "  -0.31%  " },
    @{ Cell = "D25"; Value = "159.77" },
    @{ Cell = "E25"; Value = "  -0.17%  " },
    @{ Cell = "D26"; Value = "16.63" },
    @{ Cell = "E26"; Value = "  +1.30%  " },
    @{ Cell = "D27"; Value = "7.14" },
    @{ Cell = "E27"; Value = "  +0.78%  " },
    @{ Cell = "E28"; Value = "  -0.92%  " },
    @{ Cell = "E29"; Value = "  -0.05%  " },
    @{ Cell = "D30"; Value = "0.0528" },
    @{ Cell = "E30"; Value = "  +2.22%  " },
    @{ Cell = "E31"; Value = "  -0.43%  " },
    @{ Cell = "E32"; Value = "  +1.76%  " },
    @{ Cell = "E33"; Value = "  -0.74%  " },
    @{ Cell = "E34"; Value = "  -1.27%  " },
    @{ Cell = "D35"; Value = "1.418.44" },
    @{ Cell = "E35"; Value = "  -1.26%  " },
    @{ Cell = "D36"; Value = "0.654" },
    @{ Cell = "E36"; Value = "  +2.39%  " },
    @{ Cell = "E37"; Value = "  +0.48%  " },
    @{ Cell = "E38"; Value = "  -1.44%  " },
    @{ Cell = "D39"; Value = "0.948" },
    @{ Cell = "E39"; Value = "  +3.77%  " },
    @{ Cell = "D40"; Value = "80.70" },
    @{ Cell = "E40"; Value = "  -4.30%  " },
    @{ Cell = "B41"; Value = "HuobiToken" },
    @{ Cell = "C41"; Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht" },
    @{ Cell = "D41"; Value = "2.36" },
    @{ Cell = "E41"; Value = "  +0.55%  " },
    @{ Cell = "B42"; Value = "MXToken" },
    @{ Cell = "C42"; Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx" },
    @{ Cell = "D42"; Value = "2.73" },
    @{ Cell = "E42"; Value = "  -2.37%  " },
    @{ Cell = "E43"; Value = "  +3.92%  " },
    @{ Cell = "D44"; Value = "5.96" },
    @{ Cell = "E44"; Value = "  -0.56%  " },
    @{ Cell = "D45"; Value = "108.34" },
    @{ Cell = "E45"; Value = "  +4.12%  " },
    @{ Cell = "D46"; Value = "0.0496" },
    @{ Cell = "E46"; Value = "  +0.62%  " },
    @{ Cell = "B47"; Value = "WEMIXToken" },
    @{ Cell = "C47"; Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix" },
    @{ Cell = "D47"; Value = "1.05" },
    @{ Cell = "E47"; Value = "  -1.16%  " },
    @{ Cell = "B48"; Value = "RocketPoolETH" },
    @{ Cell = "C48"; Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth" },
    @{ Cell = "D48"; Value = "1.964.10" },
    @{ Cell = "E48"; Value = "  +0.84%  " },
    @{ Cell = "E49"; Value = "  +0.36%  " },
    @{ Cell = "E51"; Value = "  +3.54%  " }
)

foreach ($u in $updates) {
    Set-CellText $ws.Range($u.Cell) $u.Value
}
